$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.018.36'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.596.87'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.80'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.480'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.30'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('D12').Value = '1.819.61'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').Value = '1.590.71'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('E15').Value = '  +2.07%  '
$ws.Range('D16').Value = '26.019.53'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '203.18'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +13.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '143.79'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -7.18%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.19'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.13'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('D36').Value = '1.129.76'
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0165'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.76%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.494'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.781'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').Value = '1.730.34'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '54.21'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.406'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.01'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').Value = '0.0₇0950'
$ws.Range('E51').Value = '  -15.01%  '
